$p = $ppt.ActivePresentation

# Slide with SlideID 553 (creationId 1427618944) is "Memory Address" -
# find it by SlideID rather than assuming a fixed ordinal position.
$s = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).SlideID -eq 553) {
        $s = $p.Slides.Item($i)
    }
}

# Clone the existing "TextBox 26" shape (a plain, no-fill, word-wrapped,
# auto-fit textbox) so the new shape inherits the same body/text formatting,
# then reposition/resize/rename/retext it into the new caption box that gets
# added right after the "Speech Bubble: Oval 12" shape (last shape on slide).
$src = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "TextBox 26") {
        $src = $s.Shapes.Item($i)
    }
}

$src.Copy()
$new = $s.Shapes.Paste().Item(1)

$new.Name = "TextBox 1"
$new.Left = 512.4303937007874
$new.Top = 246.04623047244095
$new.Width = 189.5696062992126
$new.Height = 29.081259842519685
$new.TextFrame.TextRange.Text = "1 GB ~= 1 billion byte"
